# "added heart sensor & hub"
# - Compos!B2: swap the microcontroller module from the ESP32-S3-MINI-1
#   to the ESP32-S3-WROOM-1-N8R8 (the new "hub").
# - Compos!B7: fill in the IMU row with the new "lsm6dso" part (the new
#   motion/heart-adjacent sensor).
# - Make the Compos sheet the active tab/selection (was RTC before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Compos")

# New IMU component (row 7, was empty) - added before B2 so the shared
# string table order matches ("lsm6dso" then "ESP32-S3-WROOM-1-N8R8").
$ws.Range("B7").Value = "lsm6dso"

# Swap the µC module used.
$ws.Range("B2").Value = "ESP32-S3-WROOM-1-N8R8"

# Bring the Compos sheet to the front and select B3.
$ws.Activate()
$ws.Range("B3").Select()
